# Append a new blank/template data row (row 2) below the header row on the
# active sheet ("Đơn sale chính"), extending the used range from A1:T1 to
# A1:T2.
#
# Text columns (Tiền tố, Ngày thực hiện, Cơ sở, Khách hàng, Nguồn khách,
# Tên dịch vụ, Sale chính, Sale phụ, Bác sĩ 1, Bác sĩ 2, Phụ phẫu 1,
# Phụ phẫu 2) are left as empty strings, while the numeric/money columns
# (Mã dịch vụ, Đơn giá gốc, Upsale, Đơn giá, Thanh toán lần đầu, Trả sau,
# Đã thanh toán, Dư nợ) are initialised to 0, matching the report's
# row-template formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCols = @("A","C","D","E","F","G","H","J","Q","R","S","T")
$numCols  = @("B","I","K","L","M","N","O","P")

foreach ($col in $textCols) {
    $ws.Range("${col}2").Value = ""
}

foreach ($col in $numCols) {
    $ws.Range("${col}2").Value = 0
}
